# "added arguments to separateFile.C" - reflect the renamed/updated output
# file name (no more ".csv" extension) on the filesNames sheet, and update
# the active sheet/selection state to match the saved workbook view.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("filesNames")
$ws2 = $wb.Worksheets.Item("measuredFeatures")

# allBottom_7B_config1_sorted.csv -> allBottom_7B_config1_sorted
$ws1.Range("C2").Value = "allBottom_7B_config1_sorted"

# measuredFeatures sheet keeps its previous selection but is no longer the
# active/selected tab.
$ws2.Range("B9").Select()

# filesNames becomes the active sheet with a new selection.
$ws1.Activate()
$ws1.Range("C9").Select()
